$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 40724.89
$ws.Range("J87").Value = 40724.89
$ws.Range("L87").Value = 40724.89
$ws.Range("N87").Value = -43220.89

$ws.Range("H90").Value = 40724.89
$ws.Range("J90").Value = 40724.89
$ws.Range("L90").Value = 122174.67
$ws.Range("N90").Value = -134654.67

$ws.Range("H137").Value = 2032.55
$ws.Range("J137").Value = 2378
$ws.Range("L137").Value = 7134
$ws.Range("N137").Value = -12234

$ws.Range("H138").Value = 2136.3076
$ws.Range("I138").Value = 1872.6875
$ws.Range("J138").Value = 2204.3386
$ws.Range("K138").Value = 5618.0625
$ws.Range("L138").Value = 6613.0158
$ws.Range("M138").Value = -478.0625
$ws.Range("N138").Value = -16893.0158

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1663.6052
$ws.Range("I2").Value = 1612.8889
$ws.Range("J2").Value = 1788.091
$ws.Range("K2").Value = 1612.8889
$ws.Range("L2").Value = 1788.091
$ws.Range("M2").Value = -1499.8889
$ws.Range("N2").Value = -2014.091

$ws.Range("H32").Value = 5670.385
$ws.Range("I32").Value = 6032.9766
$ws.Range("J32").Value = 3938
$ws.Range("K32").Value = 6032.9766
$ws.Range("L32").Value = 3938
$ws.Range("M32").Value = -5745.9766
$ws.Range("N32").Value = -4512

$ws.Range("H88").Value = 500527
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 500527
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 500527
$ws.Range("N88").Value = -501339
$ws.Range("M88").ClearContents()

$ws.Range("H91").Value = 500527
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 500527
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 500527
$ws.Range("N91").Value = -503335
$ws.Range("M91").ClearContents()

$ws.Range("H97").Value = 1101.2759
$ws.Range("I97").Value = 1248.2273
$ws.Range("J97").Value = 639.4286
$ws.Range("K97").Value = 1248.2273
$ws.Range("L97").Value = 639.4286
$ws.Range("M97").Value = -752.2273
$ws.Range("N97").Value = -1631.4286

$ws.Range("H116").Value = 1663.6052
$ws.Range("I116").Value = 1612.8889
$ws.Range("J116").Value = 1788.091
$ws.Range("K116").Value = 1612.8889
$ws.Range("L116").Value = 1788.091
$ws.Range("M116").Value = 681.1111000000001
$ws.Range("N116").Value = -6376.091

$ws.Range("H132").Value = 11262.057
$ws.Range("I132").Value = 1657.5385
$ws.Range("J132").Value = 38017.5
$ws.Range("K132").Value = 4972.6155
$ws.Range("L132").Value = 114052.5
$ws.Range("M132").Value = -2442.6155
$ws.Range("N132").Value = -119112.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1663.6052
$ws.Range("I3").Value = 1612.8889
$ws.Range("J3").Value = 1788.091
$ws.Range("K3").Value = 1612.8889
$ws.Range("L3").Value = 1788.091
$ws.Range("M3").Value = -1498.8889
$ws.Range("N3").Value = -2016.091

$ws.Range("H134").Value = 3825.1516
$ws.Range("I134").Value = 4152.5386
$ws.Range("J134").Value = 2609.1428
$ws.Range("K134").Value = 12457.6158
$ws.Range("L134").Value = 7827.428400000001
$ws.Range("M134").Value = -9922.6158
$ws.Range("N134").Value = -12897.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 17073.785
$ws.Range("I31").Value = 45777.89
$ws.Range("J31").Value = 3477.1052
$ws.Range("K31").Value = 45777.89
$ws.Range("L31").Value = 3477.1052
$ws.Range("M31").Value = -45482.89
$ws.Range("N31").Value = -4067.1052

$ws.Range("H34").Value = 17073.785
$ws.Range("I34").Value = 45777.89
$ws.Range("J34").Value = 3477.1052
$ws.Range("K34").Value = 45777.89
$ws.Range("L34").Value = 3477.1052
$ws.Range("M34").Value = -45575.89
$ws.Range("N34").Value = -3881.1052

$ws.Range("H69").Value = 4752.4287
$ws.Range("J69").Value = 9444
$ws.Range("L69").Value = 9444
$ws.Range("N69").Value = -10942

$ws.Range("H72").Value = 4752.4287
$ws.Range("J72").Value = 9444
$ws.Range("L72").Value = 28332
$ws.Range("N72").Value = -35820

$ws.Range("H94").Value = 2435.5625
$ws.Range("J94").Value = 2872.3333
$ws.Range("L94").Value = 2872.3333
$ws.Range("N94").Value = -3774.3333

$ws.Range("H122").Value = 1170.9678
$ws.Range("I122").Value = 1046.579
$ws.Range("K122").Value = 3139.737
$ws.Range("M122").Value = -689.7370000000001

$ws.Range("H134").Value = 960.625
$ws.Range("I134").Value = 812.1429000000001
$ws.Range("J134").Value = 2000
$ws.Range("K134").Value = 2436.4287
$ws.Range("L134").Value = 6000
$ws.Range("M134").Value = 98.57129999999961
$ws.Range("N134").Value = -11070

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H49").Value = 2000
$ws.Range("J49").Value = 2000
$ws.Range("L49").Value = 6000
$ws.Range("N49").Value = -6312

$ws.Range("H113").Value = 1012.5
$ws.Range("J113").Value = 1028.5714
$ws.Range("L113").Value = 3085.7142
$ws.Range("N113").Value = -7425.7142

$ws.Range("H131").Value = 808.45
$ws.Range("I131").Value = 624.75
$ws.Range("J131").Value = 816.1042
$ws.Range("K131").Value = 1874.25
$ws.Range("L131").Value = 2448.3126
$ws.Range("M131").Value = 3165.75
$ws.Range("N131").Value = -12528.3126

$ws.Range("H139").Value = 1531.3
$ws.Range("I139").Value = 1194.5927
$ws.Range("J139").Value = 4561.6665
$ws.Range("K139").Value = 3583.7781
$ws.Range("L139").Value = 13684.9995
$ws.Range("M139").Value = 1556.2219
$ws.Range("N139").Value = -23964.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4106.75
$ws.Range("I126").Value = 3357.238
$ws.Range("J126").Value = 5537.636
$ws.Range("K126").Value = 10071.714
$ws.Range("L126").Value = 16612.908
$ws.Range("M126").Value = -7601.714
$ws.Range("N126").Value = -21552.908

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 399.4375
$ws.Range("I16").Value = 428
$ws.Range("K16").Value = 428
$ws.Range("M16").Value = -258

$ws.Range("H136").Value = 47554.453
$ws.Range("I136").Value = 84699
$ws.Range("J136").Value = 2981
$ws.Range("K136").Value = 254097
$ws.Range("L136").Value = 8943
$ws.Range("M136").Value = -251547
$ws.Range("N136").Value = -14043

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 0
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 0
$ws.Range("M96").ClearContents()
$ws.Range("N96").ClearContents()

$ws.Range("H132").Value = 1871.1538
$ws.Range("I132").Value = 1234.7
$ws.Range("J132").Value = 3992.6667
$ws.Range("K132").Value = 3704.1
$ws.Range("L132").Value = 11978.0001
$ws.Range("M132").Value = -1174.1
$ws.Range("N132").Value = -17038.0001

$ws.Range("H136").Value = 27028396
$ws.Range("J136").Value = 2066.5
$ws.Range("L136").Value = 6199.5
$ws.Range("N136").Value = -11299.5
